# Auto-generated Excel COM-interop script applying the Ultima_Profits price/profit
# recalculation update across the ALC, ARM, BSM, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 337436.06
$ws.Range("I127").Value = 494.5
$ws.Range("K127").Value = 1483.5
$ws.Range("M127").Value = 3476.5

$ws.Range("H128").Value = 70520
$ws.Range("J128").Value = 70520
$ws.Range("L128").Value = 70520
$ws.Range("N128").Value = -80480

$ws.Range("H129").Value = 2753.3333
$ws.Range("I129").Value = 650
$ws.Range("J129").Value = 3076.923
$ws.Range("K129").Value = 1950
$ws.Range("L129").Value = 9230.769
$ws.Range("M129").Value = 3050
$ws.Range("N129").Value = -19230.769

$ws.Range("H131").Value = 5971.697
$ws.Range("I131").Value = 884.125
$ws.Range("J131").Value = 10760
$ws.Range("K131").Value = 2652.375
$ws.Range("L131").Value = 32280
$ws.Range("M131").Value = 2387.625
$ws.Range("N131").Value = -42360

$ws.Range("H132").Value = 4810778
$ws.Range("I132").Value = 2683.8262
$ws.Range("J132").Value = 41672830
$ws.Range("K132").Value = 8051.4786
$ws.Range("L132").Value = 125018490
$ws.Range("M132").Value = -5521.4786
$ws.Range("N132").Value = -125023550

$ws.Range("H135").Value = 1005.72546
$ws.Range("J135").Value = 1038.8
$ws.Range("L135").Value = 9349.199999999999
$ws.Range("N135").Value = -14419.2

$ws.Range("H137").Value = 9551.348
$ws.Range("I137").Value = 936.9048
$ws.Range("J137").Value = 100003
$ws.Range("K137").Value = 2810.7144
$ws.Range("L137").Value = 300009
$ws.Range("M137").Value = -260.7143999999998
$ws.Range("N137").Value = -305109

$ws.Range("H138").Value = 4571027.5
$ws.Range("I138").Value = 15154676
$ws.Range("J138").Value = 5532.137
$ws.Range("K138").Value = 45464028
$ws.Range("L138").Value = 16596.411
$ws.Range("M138").Value = -45458888
$ws.Range("N138").Value = -26876.411

$ws.Range("H140").Value = 52857.145
$ws.Range("J140").Value = 59166.668
$ws.Range("L140").Value = 59166.668
$ws.Range("N140").Value = -69526.66800000001

$ws.Range("H141").Value = 1222.9474
$ws.Range("I141").Value = 878.44116
$ws.Range("K141").Value = 2635.32348
$ws.Range("M141").Value = 2544.67652

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1511.8727
$ws.Range("I132").Value = 1182.8
$ws.Range("K132").Value = 3548.4
$ws.Range("M132").Value = -1018.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 7875
$ws.Range("J118").Value = 7875
$ws.Range("L118").Value = 7875
$ws.Range("N118").Value = -11189

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 15986.667
$ws.Range("J95").Value = 15986.667
$ws.Range("L95").Value = 47960.001
$ws.Range("N95").Value = -52078.001

$ws.Range("H114").Value = 1552.619
$ws.Range("I114").Value = 1286.7
$ws.Range("J114").Value = 1794.3636
$ws.Range("K114").Value = 3860.1
$ws.Range("L114").Value = 5383.0908
$ws.Range("M114").Value = -606.1000000000004
$ws.Range("N114").Value = -11891.0908

$ws.Range("H117").Value = 385.45456
$ws.Range("I117").Value = 328
$ws.Range("J117").Value = 433.33334
$ws.Range("K117").Value = 984
$ws.Range("L117").Value = 1300.00002
$ws.Range("M117").Value = 2458
$ws.Range("N117").Value = -8184.000019999999

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -39880
$ws.Range("M126").ClearContents()

$ws.Range("H129").Value = 5094.609
$ws.Range("I129").Value = 6801.25
$ws.Range("J129").Value = 4184.4
$ws.Range("K129").Value = 20403.75
$ws.Range("L129").Value = 12553.2
$ws.Range("M129").Value = -15403.75
$ws.Range("N129").Value = -22553.2

$ws.Range("H134").Value = 5593.625
$ws.Range("I134").Value = 2753.9167
$ws.Range("J134").Value = 8433.333000000001
$ws.Range("K134").Value = 8261.750100000001
$ws.Range("L134").Value = 25299.999
$ws.Range("M134").Value = -3191.750100000001
$ws.Range("N134").Value = -35439.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7772.143
$ws.Range("J57").Value = 16000
$ws.Range("L57").Value = 16000
$ws.Range("N57").Value = -17640

$ws.Range("H70").Value = 9656.6
$ws.Range("I70").Value = 12610.087
$ws.Range("J70").Value = 3995.75
$ws.Range("K70").Value = 12610.087
$ws.Range("L70").Value = 3995.75
$ws.Range("M70").Value = -12340.087
$ws.Range("N70").Value = -4535.75

$ws.Range("H73").Value = 9656.6
$ws.Range("I73").Value = 12610.087
$ws.Range("J73").Value = 3995.75
$ws.Range("K73").Value = 12610.087
$ws.Range("L73").Value = 3995.75
$ws.Range("M73").Value = -11674.087
$ws.Range("N73").Value = -5867.75

$ws.Range("H80").Value = 2906.8
$ws.Range("I80").Value = 2712.875
$ws.Range("J80").Value = 3128.4285
$ws.Range("K80").Value = 2712.875
$ws.Range("L80").Value = 3128.4285
$ws.Range("M80").Value = -1714.875
$ws.Range("N80").Value = -5124.4285

$ws.Range("H83").Value = 2906.8
$ws.Range("I83").Value = 2712.875
$ws.Range("J83").Value = 3128.4285
$ws.Range("K83").Value = 13564.375
$ws.Range("L83").Value = 15642.1425
$ws.Range("M83").Value = -8572.375
$ws.Range("N83").Value = -25626.1425

$ws.Range("H126").Value = 4218.1875
$ws.Range("I126").Value = 3174.5
$ws.Range("J126").Value = 4566.0835
$ws.Range("K126").Value = 9523.5
$ws.Range("L126").Value = 13698.2505
$ws.Range("M126").Value = -7053.5
$ws.Range("N126").Value = -18638.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8694.049999999999
$ws.Range("J7").Value = 12237.5
$ws.Range("L7").Value = 12237.5
$ws.Range("N7").Value = -12461.5

$ws.Range("H35").Value = 1900
$ws.Range("I35").Value = 1900
$ws.Range("J35").Value = 1900
$ws.Range("K35").Value = 1900
$ws.Range("L35").Value = 1900
$ws.Range("M35").Value = -1564
$ws.Range("N35").Value = -2572

$ws.Range("H40").Value = 4994.25
$ws.Range("I40").Value = 14050
$ws.Range("J40").Value = 1975.6666
$ws.Range("K40").Value = 14050
$ws.Range("L40").Value = 1975.6666
$ws.Range("M40").Value = -13914
$ws.Range("N40").Value = -2247.6666

$ws.Range("H55").Value = 304.25925
$ws.Range("I55").Value = 175.8
$ws.Range("J55").Value = 464.83334
$ws.Range("K55").Value = 175.8
$ws.Range("L55").Value = 464.83334
$ws.Range("M55").Value = -2.800000000000011
$ws.Range("N55").Value = -810.83334

$ws.Range("H126").Value = 8694.049999999999
$ws.Range("J126").Value = 12237.5
$ws.Range("L126").Value = 36712.5
$ws.Range("N126").Value = -41652.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 5007
$ws.Range("J15").Value = 5007
$ws.Range("L15").Value = 5007
$ws.Range("N15").Value = -5583

$ws.Range("H54").Value = 21035.715
$ws.Range("I54").Value = 9753.333000000001
$ws.Range("J54").Value = 29497.5
$ws.Range("K54").Value = 9753.333000000001
$ws.Range("L54").Value = 29497.5
$ws.Range("M54").Value = -9233.333000000001
$ws.Range("N54").Value = -30537.5

$ws.Range("H81").Value = 940
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 2700
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 5400
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -7522

$ws.Range("H84").Value = 940
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 2700
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 27000
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -37608

$ws.Range("H136").Value = 1158.4036
$ws.Range("I136").Value = 784.88464
$ws.Range("J136").Value = 5043
$ws.Range("K136").Value = 2354.65392
$ws.Range("L136").Value = 15129
$ws.Range("M136").Value = 195.3460800000003
$ws.Range("N136").Value = -20229
